$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated dSF (column F) values per row, as redetermined after repulling data
$values = @{
    2  = -2
    3  = 0
    5  = -2
    6  = 1
    7  = 4
    8  = 3
    9  = -1
    10 = -1
    11 = 1
    12 = -3
    13 = -1
    14 = 5
    15 = 1
    16 = 3
    17 = -5
    18 = -2
    19 = 4
    20 = 1
    22 = 3
    23 = -4
    24 = 2
    25 = 5
    26 = 7
    27 = 1
    28 = 3
    29 = 2
    30 = 2
    31 = 1
    32 = -3
    33 = 2
    34 = 6
    35 = 1
    38 = 2
    40 = 3
}

foreach ($row in $values.Keys) {
    $ws.Range("F$row").Value = $values[$row]
}
